# edit.ps1 - small fixes in presentation
#
#  1) Five occurrences of the code token "connect()" lose their empty
#     call-parentheses and become plain "connect" (slides 32, 36 (x2),
#     38, 40). A sixth, decoy occurrence of "connect()" immediately
#     followed by "(App)" on slide 32 is left untouched.
#  2) On slide 5, the second of three "SingleName" mentions (the one
#     talking about the *array* example, "массив в SingleName") was a
#     copy/paste slip - it should reference the "MultiNames" component
#     instead, since that is the component that actually holds the
#     array of names.

$p = $ppt.ActivePresentation

function Replace-StandaloneConnect {
    param(
        [int]$SlideIndex,
        [int]$MaxReplacements
    )

    $slide = $p.Slides.Item($SlideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if (-not $shape.HasTextFrame) { continue }

        $tr = $shape.TextFrame.TextRange
        $text = $tr.Text
        if ($text.IndexOf("connect()") -lt 0) { continue }

        $done = 0
        $searchFrom = 0
        while ($done -lt $MaxReplacements) {
            $idx = $text.IndexOf("connect()", $searchFrom)
            if ($idx -lt 0) { break }

            $afterPos = $idx + 9
            $followedByParen = ($afterPos -lt $text.Length) -and ($text.Substring($afterPos, 1) -eq "(")

            if ($followedByParen) {
                # skip the "connect()(App)" decoy - keep scanning past it
                $searchFrom = $idx + 1
                continue
            }

            # Characters() is 1-based; replace just the "()" worth of
            # trailing characters so the run's formatting (lang/dirty/err)
            # is preserved and only the text content changes.
            $tr.Characters($idx + 1, 9).Text = "connect"

            # text got 2 characters shorter - refresh and keep scanning
            # forward from the same spot in case of further matches
            $text = $tr.Text
            $searchFrom = $idx + 7
            $done = $done + 1
        }
    }
}

# slide 32: "При использовании connect() у нас не было повторной..."
Replace-StandaloneConnect 32 1

# slide 36: two mentions, both "..., обернутого в connect() ..."
Replace-StandaloneConnect 36 2

# slide 38: "...и удаления обертки connect() мы видим..."
Replace-StandaloneConnect 38 1

# slide 40: "...одновременном использовании обертки connect() и memo/PureComponent"
Replace-StandaloneConnect 40 1

# slide 5: fix the "массив в SingleName" -> "массив в MultiNames" slip.
# There are three "SingleName" runs in this shape; only the middle one
# (preceded by "массив в ") refers to the array example and must change.
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shape = $slide5.Shapes.Item($i)
    if (-not $shape.HasTextFrame) { continue }

    $tr = $shape.TextFrame.TextRange
    $text = $tr.Text

    $idx1 = $text.IndexOf("SingleName")
    if ($idx1 -lt 0) { continue }
    $idx2 = $text.IndexOf("SingleName", $idx1 + 1)
    if ($idx2 -lt 0) { continue }

    $tr.Characters($idx2 + 1, 10).Text = "MultiNames"
    break
}
